$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '64.292.33'
$ws.Range('E2').Value = '  -4.28%  '
Set-TextValue $ws.Range('D3') '3.017.79'
$ws.Range('E3').Value = '  -6.91%  '
$ws.Range('E4').Value = '  +0.31%  '
Set-TextValue $ws.Range('D5') '554.64'
$ws.Range('E5').Value = '  -6.99%  '
Set-TextValue $ws.Range('D6') '139.67'
$ws.Range('E6').Value = '  -9.09%  '
Set-TextValue $ws.Range('D7') '1.00'
$ws.Range('E7').Value = '  +0.14%  '
Set-TextValue $ws.Range('D8') '3.003.42'
$ws.Range('E8').Value = '  -7.16%  '
$ws.Range('E9').Value = '  -12.20%  '
$ws.Range('E10').Value = '  -11.89%  '
Set-TextValue $ws.Range('D11') '5.96'
$ws.Range('E11').Value = '  -10.82%  '
Set-TextValue $ws.Range('D12') '0.448'
$ws.Range('E12').Value = '  -11.08%  '
Set-TextValue $ws.Range('D13') '34.14'
$ws.Range('E13').Value = '  -13.38%  '
Set-TextValue $ws.Range('D14') '0.0000213'
$ws.Range('E14').Value = '  -14.05%  '
Set-TextValue $ws.Range('D15') '3.515.30'
$ws.Range('E15').Value = '  -6.40%  '
Set-TextValue $ws.Range('D16') '64.383.64'
$ws.Range('E16').Value = '  -4.08%  '
$ws.Range('E17').Value = '  -4.00%  '
Set-TextValue $ws.Range('D18') '3.031.03'
$ws.Range('E18').Value = '  -6.51%  '
Set-TextValue $ws.Range('D19') '481.76'
$ws.Range('E19').Value = '  -9.81%  '
Set-TextValue $ws.Range('D20') '6.45'
$ws.Range('E20').Value = '  -10.90%  '
Set-TextValue $ws.Range('D21') '13.38'
$ws.Range('E21').Value = '  -11.85%  '
Set-TextValue $ws.Range('D22') '0.656'
$ws.Range('E22').Value = '  -14.41%  '
Set-TextValue $ws.Range('D23') '6.93'
$ws.Range('E23').Value = '  -12.52%  '
Set-TextValue $ws.Range('D24') '12.42'
$ws.Range('E24').Value = '  -10.94%  '
Set-TextValue $ws.Range('D25') '77.71'
$ws.Range('E25').Value = '  -9.85%  '
Set-TextValue $ws.Range('D26') '0.997'
$ws.Range('E26').Value = '  -0.14%  '
Set-TextValue $ws.Range('D27') '2.70'
$ws.Range('E27').Value = '  -15.43%  '
$ws.Range('E28').Value = '  -7.18%  '
Set-TextValue $ws.Range('D29') '7.54'
$ws.Range('E29').Value = '  -8.49%  '
Set-TextValue $ws.Range('D30') '25.57'
$ws.Range('E30').Value = '  -13.43%  '
Set-TextValue $ws.Range('D31') '2.55'
$ws.Range('E31').Value = '  -3.56%  '
Set-TextValue $ws.Range('D32') '1.10'
$ws.Range('E32').Value = '  -4.64%  '
$ws.Range('E33').Value = '  -0.16%  '
Set-TextValue $ws.Range('D34') '509.22'
$ws.Range('E34').Value = '  -5.76%  '
Set-TextValue $ws.Range('D35') '5.24'
$ws.Range('E35').Value = '  -9.52%  '
Set-TextValue $ws.Range('D36') '51.82'
$ws.Range('E36').Value = '  -2.79%  '
Set-TextValue $ws.Range('D37') '5.74'
$ws.Range('E37').Value = '  -13.13%  '
$ws.Range('E38').Value = '  -5.84%  '
Set-TextValue $ws.Range('D39') '0.0791'
$ws.Range('E39').Value = '  -9.85%  '
Set-TextValue $ws.Range('D40') '0.119'
$ws.Range('E40').Value = '  -7.34%  '
Set-TextValue $ws.Range('D41') '8.09'
$ws.Range('E41').Value = '  -13.74%  '
Set-TextValue $ws.Range('D42') '2.802.72'
$ws.Range('E42').Value = '  -5.01%  '
$ws.Range('E43').Value = '  -11.16%  '
$ws.Range('E44').Value = '  -0.14%  '
$ws.Range('E45').Value = '  -12.12%  '
Set-TextValue $ws.Range('D46') '1.98'
$ws.Range('E46').Value = '  -7.66%  '
$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range('D47') '0.0₃0509'
$ws.Range('E47').Value = '  -14.64%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D48') '0.106'
$ws.Range('E48').Value = '  -8.13%  '
Set-TextValue $ws.Range('D49') '115.94'
$ws.Range('E49').Value = '  -5.78%  '
Set-TextValue $ws.Range('D50') '23.47'
$ws.Range('E50').Value = '  -12.17%  '
$ws.Range('E51').Value = '  -17.95%  '
